# executive-presentation.pptx edit
#
# Final state: the deck is trimmed down to 4 slides:
#   1. Title slide (unchanged)
#   2. "Why This Solution?" comparison table (previously slide 7)
#   3. "Business Value - Financial Impact" table (previously slide 8)
#   4. "Risk Mitigation" table (previously slide 11)
# All other slides (2-6, 9-10, 12-17) are removed, and on the three
# surviving table slides the bold formatting is stripped from the last
# data row (and, for slides 3/4, the header row too).

$p = $ppt.ActivePresentation

# --- 1. Remove the slides we don't want, from the highest index down so
#        that the indices of slides we still need to delete don't shift. ---
$slidesToDelete = @(17, 16, 15, 14, 13, 12, 10, 9, 6, 5, 4, 3, 2)
foreach ($idx in $slidesToDelete) {
    $p.Slides.Item($idx).Delete()
}

# After the deletions above, the four remaining slides (originally
# 1, 7, 8, 11) now sit at positions 1, 2, 3, 4 respectively.

# --- 2. Slide 2 ("Why This Solution?"): un-bold the last comparison row ---
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
for ($c = 1; $c -le $tbl2.Columns.Count; $c++) {
    $tbl2.Cell(4, $c).Shape.TextFrame.TextRange.Font.Bold = $false
}

# --- 3. Slide 3 ("Business Value - Financial Impact"): un-bold header row
#        and the final ROI row ---
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
for ($c = 1; $c -le $tbl3.Columns.Count; $c++) {
    $tbl3.Cell(1, $c).Shape.TextFrame.TextRange.Font.Bold = $false
    $tbl3.Cell(6, $c).Shape.TextFrame.TextRange.Font.Bold = $false
}

# --- 4. Slide 4 ("Risk Mitigation"): un-bold header row and the final
#        [Risk 3] row ---
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
for ($c = 1; $c -le $tbl4.Columns.Count; $c++) {
    $tbl4.Cell(1, $c).Shape.TextFrame.TextRange.Font.Bold = $false
    $tbl4.Cell(4, $c).Shape.TextFrame.TextRange.Font.Bold = $false
}

Write-Host "Final slide count: $($p.Slides.Count)"
